$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's data (2026/01/30, Friday) was inserted right after the last
# 2026/01/29 row, pushing every following row down by one (721 -> 722 ... 762 -> 763).
$ws.Rows.Item(721).Insert()

# Column A holds dates stored as plain text (e.g. "2026/12/29"), not real
# Excel dates. Assigning a literal string like "2026/01/30" to .Value would
# make Excel auto-coerce it into a date serial number/date-formatted cell,
# so instead stage the text on a scratch cell via a formula that evaluates
# to a text string (no auto date-coercion, no quote-prefix style needed),
# copy the computed value over with Copy/PasteSpecial, and then discard the
# scratch row.
$scratchRow = 1000
$ws.Cells.Item($scratchRow, 1).Formula = "=""2026/01/30"""
$ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, 1)).Copy()
$ws.Range("A721").PasteSpecial()
$ws.Rows.Item($scratchRow).Delete()

$ws.Cells.Item(721, 2).Value = "金"
$ws.Cells.Item(721, 3).Value = 2
$ws.Cells.Item(721, 4).Value = 21
